# Insert a new weekly price-report row for "Perejil" (Vega Central Mapocho de
# Santiago) just before the existing row 448. This shifts every following
# record down by one row (old row 448 -> 449, ..., old row 536 -> 537) and
# grows the sheet's used range from A1:R536 to A1:R537.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a blank row at position 448 (pushes rows 448..536 down to 449..537).
$ws.Rows.Item(448).Insert()

# Populate the newly inserted row 448 with the new weekly record.
$ws.Cells.Item(448, 1).Value = 9
$ws.Cells.Item(448, 2).Value = "Vega Central Mapocho de Santiago"
$ws.Cells.Item(448, 3).Value = "Metropolitana"
$ws.Cells.Item(448, 4).Value = 45015
$ws.Cells.Item(448, 5).Value = 13
$ws.Cells.Item(448, 6).Value = 100112044
$ws.Cells.Item(448, 7).Value = "Perejil"
$ws.Cells.Item(448, 8).Value = "Sin especificar"
$ws.Cells.Item(448, 9).Value = "Primera"
$ws.Cells.Item(448, 10).Value = 70
$ws.Cells.Item(448, 11).Value = 13000
$ws.Cells.Item(448, 12).Value = 14000
$ws.Cells.Item(448, 13).Value = 13500
$ws.Cells.Item(448, 14).Value = "$/docena de atados"
$ws.Cells.Item(448, 15).Value = "Región Metropolitana"
$ws.Cells.Item(448, 16).Value = 4500
$ws.Cells.Item(448, 17).Value = 3
$ws.Cells.Item(448, 18).Value = "Hortaliza"
